{"js": "// Shrink the \"Heading 2\" / \"Heading 3\" report styles (and their linked\n// character styles) down to a uniform 14pt (w:sz 28 half-points), matching\n// the updated Plasmix report template:\n//   heading 2 / \"\u6807\u9898 2 \u5b57\u7b26\" : 18pt (36) -> 14pt (28)\n//   heading 3 / \"\u6807\u9898 3 \u5b57\u7b26\" : 16pt (32) -> 14pt (28)\n\nconst styles = context.document.getStyles();\n\nconst heading2 = styles.getByNameOrNullObject(\"Heading 2\");\nconst heading3 = styles.getByNameOrNullObject(\"Heading 3\");\nconst heading2Char = styles.getByNameOrNullObject(\"\u6807\u9898 2 \u5b57\u7b26\");\nconst heading3Char = styles.getByNameOrNullObject(\"\u6807\u9898 3 \u5b57\u7b26\");\n\nheading2.load(\"isNullObject\");\nheading3.load(\"isNullObject\");\nheading2Char.load(\"isNullObject\");\nheading3Char.load(\"isNullObject\");\nawait context.sync();\n\nconst NEW_SIZE = 14; // points == w:sz 28 (half-points)\n\nif (!heading2.isNullObject) {\n  heading2.font.size = NEW_SIZE;\n}\nif (!heading3.isNullObject) {\n  heading3.font.size = NEW_SIZE;\n}\nif (!heading2Char.isNullObject) {\n  heading2Char.font.size = NEW_SIZE;\n}\nif (!heading3Char.isNullObject) {\n  heading3Char.font.size = NEW_SIZE;\n}\n\nawait context.sync();\n", "ps1": "# Shrink the \"Heading 2\" / \"Heading 3\" report styles (and their linked\n# character styles) down to a uniform 14pt (w:sz 28 half-points), matching\n# the updated Plasmix report template:\n#   heading 2 / \"\u6807\u9898 2 \u5b57\u7b26\" : 18pt (36) -> 14pt (28)\n#   heading 3 / \"\u6807\u9898 3 \u5b57\u7b26\" : 16pt (32) -> 14pt (28)\n\n$d = $word.ActiveDocument\n\n$NewSize = 14  # points == w:sz 28 (half-points)\n\n# Paragraph style \"Heading 2\"\n$heading2 = $d.Styles(\"Heading 2\")\nif ($heading2 -ne $null) {\n    $heading2.Font.Size = $NewSize\n}\n\n# Paragraph style \"Heading 3\"\n$heading3 = $d.Styles(\"Heading 3\")\nif ($heading3 -ne $null) {\n    $heading3.Font.Size = $NewSize\n}\n\n# Linked character style for Heading 2 (\"\u6807\u9898 2 \u5b57\u7b26\" = \"Heading 2 Char\")\n$heading2Char = $d.Styles(\"\u6807\u9898 2 \u5b57\u7b26\")\nif ($heading2Char -ne $null) {\n    $heading2Char.Font.Size = $NewSize\n    $heading2Char.Font.LanguageIDFarEast = \"en-US\"\n}\n\n# Linked character style for Heading 3 (\"\u6807\u9898 3 \u5b57\u7b26\" = \"Heading 3 Char\")\n$heading3Char = $d.Styles(\"\u6807\u9898 3 \u5b57\u7b26\")\nif ($heading3Char -ne $null) {\n    $heading3Char.Font.Size = $NewSize\n    $heading3Char.Font.SizeBi = 12\n    $heading3Char.Font.LanguageIDFarEast = \"en-US\"\n}\n"}
